$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value2 = "ECs"
$ws.Cells.Item(2, 2).Value2 = "Clec11a"
$ws.Cells.Item(2, 3).Value2 = "Itga10"
$ws.Cells.Item(2, 4).Value2 = "ECs"
$ws.Cells.Item(2, 5).Value2 = 2
$ws.Cells.Item(2, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(2, 7).Value2 = 0.152918
$ws.Cells.Item(2, 8).Value2 = 0.458754
$ws.Cells.Item(2, 9).Value2 = 0.007158013163202275
$ws.Cells.Item(2, 10).Value2 = 0.007158013163202275
$ws.Cells.Item(2, 11).Value2 = 3
$ws.Cells.Item(2, 12).Value2 = 1
$ws.Cells.Item(2, 13).Value2 = 0.970007
$ws.Cells.Item(2, 14).Value2 = 2.910021
$ws.Cells.Item(2, 15).Value2 = 0.1871949781160357
$ws.Cells.Item(2, 16).Value2 = 0.1871949781160357
$ws.Cells.Item(2, 17).Value2 = 0.148331530426
$ws.Cells.Item(2, 18).Value2 = 1.334983773834
$ws.Cells.Item(2, 19).Value2 = 0.001339944117439946
$ws.Cells.Item(2, 20).Value2 = 0.001339944117439946

# Row 3
$ws.Cells.Item(3, 1).Value2 = "ECs"
$ws.Cells.Item(3, 2).Value2 = "Clec11a"
$ws.Cells.Item(3, 3).Value2 = "Itga10"
$ws.Cells.Item(3, 4).Value2 = "FAPs"
$ws.Cells.Item(3, 5).Value2 = 2
$ws.Cells.Item(3, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(3, 7).Value2 = 0.152918
$ws.Cells.Item(3, 8).Value2 = 0.458754
$ws.Cells.Item(3, 9).Value2 = 0.007158013163202275
$ws.Cells.Item(3, 10).Value2 = 0.007158013163202275
$ws.Cells.Item(3, 11).Value2 = 3
$ws.Cells.Item(3, 12).Value2 = 1
$ws.Cells.Item(3, 13).Value2 = 2.313512
$ws.Cells.Item(3, 14).Value2 = 6.940536
$ws.Cells.Item(3, 15).Value2 = 0.4464687659070358
$ws.Cells.Item(3, 16).Value2 = 0.4464687659070357
$ws.Cells.Item(3, 17).Value2 = 0.353777628016
$ws.Cells.Item(3, 18).Value2 = 3.183998652144
$ws.Cells.Item(3, 19).Value2 = 0.003195829303321237
$ws.Cells.Item(3, 20).Value2 = 0.003195829303321237

# Row 4
$ws.Cells.Item(4, 1).Value2 = "ECs"
$ws.Cells.Item(4, 2).Value2 = "Clec11a"
$ws.Cells.Item(4, 3).Value2 = "Itga10"
$ws.Cells.Item(4, 4).Value2 = "MuSCs"
$ws.Cells.Item(4, 5).Value2 = 2
$ws.Cells.Item(4, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(4, 7).Value2 = 0.152918
$ws.Cells.Item(4, 8).Value2 = 0.458754
$ws.Cells.Item(4, 9).Value2 = 0.007158013163202275
$ws.Cells.Item(4, 10).Value2 = 0.007158013163202275
$ws.Cells.Item(4, 11).Value2 = 3
$ws.Cells.Item(4, 12).Value2 = 1
$ws.Cells.Item(4, 13).Value2 = 1.616694666666667
$ws.Cells.Item(4, 14).Value2 = 4.850084
$ws.Cells.Item(4, 15).Value2 = 0.3119947822510337
$ws.Cells.Item(4, 16).Value2 = 0.3119947822510336
$ws.Cells.Item(4, 17).Value2 = 0.2472217150373333
$ws.Cells.Item(4, 18).Value2 = 2.224995435336
$ws.Cells.Item(4, 19).Value2 = 0.002233262758203327
$ws.Cells.Item(4, 20).Value2 = 0.002233262758203326

# Row 5
$ws.Cells.Item(5, 1).Value2 = "ECs"
$ws.Cells.Item(5, 2).Value2 = "Clec11a"
$ws.Cells.Item(5, 3).Value2 = "Itga10"
$ws.Cells.Item(5, 4).Value2 = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value2 = 2
$ws.Cells.Item(5, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(5, 7).Value2 = 0.152918
$ws.Cells.Item(5, 8).Value2 = 0.458754
$ws.Cells.Item(5, 9).Value2 = 0.007158013163202275
$ws.Cells.Item(5, 10).Value2 = 0.007158013163202275
$ws.Cells.Item(5, 11).Value2 = 3
$ws.Cells.Item(5, 12).Value2 = 1
$ws.Cells.Item(5, 13).Value2 = 0.2815866666666667
$ws.Cells.Item(5, 14).Value2 = 0.8447600000000001
$ws.Cells.Item(5, 15).Value2 = 0.05434147372589489
$ws.Cells.Item(5, 16).Value2 = 0.05434147372589488
$ws.Cells.Item(5, 17).Value2 = 0.04305966989333334
$ws.Cells.Item(5, 18).Value2 = 0.38753702904
$ws.Cells.Item(5, 19).Value2 = 0.0003889769842377662
$ws.Cells.Item(5, 20).Value2 = 0.0003889769842377661

# Row 6
$ws.Cells.Item(6, 1).Value2 = "FAPs"
$ws.Cells.Item(6, 2).Value2 = "Clec11a"
$ws.Cells.Item(6, 3).Value2 = "Itga10"
$ws.Cells.Item(6, 4).Value2 = "ECs"
$ws.Cells.Item(6, 5).Value2 = 3
$ws.Cells.Item(6, 6).Value2 = 1
$ws.Cells.Item(6, 7).Value2 = 20.82581466666667
$ws.Cells.Item(6, 8).Value2 = 62.47744400000001
$ws.Cells.Item(6, 9).Value2 = 0.9748457050079848
$ws.Cells.Item(6, 10).Value2 = 0.9748457050079848
$ws.Cells.Item(6, 11).Value2 = 3
$ws.Cells.Item(6, 12).Value2 = 1
$ws.Cells.Item(6, 13).Value2 = 0.970007
$ws.Cells.Item(6, 14).Value2 = 2.910021
$ws.Cells.Item(6, 15).Value2 = 0.1871949781160357
$ws.Cells.Item(6, 16).Value2 = 0.1871949781160357
$ws.Cells.Item(6, 17).Value2 = 20.20118600736934
$ws.Cells.Item(6, 18).Value2 = 181.810674066324
$ws.Cells.Item(6, 19).Value2 = 0.1824862204154812
$ws.Cells.Item(6, 20).Value2 = 0.1824862204154811

# Row 7
$ws.Cells.Item(7, 1).Value2 = "FAPs"
$ws.Cells.Item(7, 2).Value2 = "Clec11a"
$ws.Cells.Item(7, 3).Value2 = "Itga10"
$ws.Cells.Item(7, 4).Value2 = "FAPs"
$ws.Cells.Item(7, 5).Value2 = 3
$ws.Cells.Item(7, 6).Value2 = 1
$ws.Cells.Item(7, 7).Value2 = 20.82581466666667
$ws.Cells.Item(7, 8).Value2 = 62.47744400000001
$ws.Cells.Item(7, 9).Value2 = 0.9748457050079848
$ws.Cells.Item(7, 10).Value2 = 0.9748457050079848
$ws.Cells.Item(7, 11).Value2 = 3
$ws.Cells.Item(7, 12).Value2 = 1
$ws.Cells.Item(7, 13).Value2 = 2.313512
$ws.Cells.Item(7, 14).Value2 = 6.940536
$ws.Cells.Item(7, 15).Value2 = 0.4464687659070358
$ws.Cells.Item(7, 16).Value2 = 0.4464687659070357
$ws.Cells.Item(7, 17).Value2 = 48.18077214110934
$ws.Cells.Item(7, 18).Value2 = 433.626949269984
$ws.Cells.Item(7, 19).Value2 = 0.4352381588646893
$ws.Cells.Item(7, 20).Value2 = 0.4352381588646892

# Row 8
$ws.Cells.Item(8, 1).Value2 = "FAPs"
$ws.Cells.Item(8, 2).Value2 = "Clec11a"
$ws.Cells.Item(8, 3).Value2 = "Itga10"
$ws.Cells.Item(8, 4).Value2 = "MuSCs"
$ws.Cells.Item(8, 5).Value2 = 3
$ws.Cells.Item(8, 6).Value2 = 1
$ws.Cells.Item(8, 7).Value2 = 20.82581466666667
$ws.Cells.Item(8, 8).Value2 = 62.47744400000001
$ws.Cells.Item(8, 9).Value2 = 0.9748457050079848
$ws.Cells.Item(8, 10).Value2 = 0.9748457050079848
$ws.Cells.Item(8, 11).Value2 = 3
$ws.Cells.Item(8, 12).Value2 = 1
$ws.Cells.Item(8, 13).Value2 = 1.616694666666667
$ws.Cells.Item(8, 14).Value2 = 4.850084
$ws.Cells.Item(8, 15).Value2 = 0.3119947822510337
$ws.Cells.Item(8, 16).Value2 = 0.3119947822510336
$ws.Cells.Item(8, 17).Value2 = 33.66898350058845
$ws.Cells.Item(8, 18).Value2 = 303.020851505296
$ws.Cells.Item(8, 19).Value2 = 0.3041467734623216
$ws.Cells.Item(8, 20).Value2 = 0.3041467734623215

# Row 9
$ws.Cells.Item(9, 1).Value2 = "FAPs"
$ws.Cells.Item(9, 2).Value2 = "Clec11a"
$ws.Cells.Item(9, 3).Value2 = "Itga10"
$ws.Cells.Item(9, 4).Value2 = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value2 = 3
$ws.Cells.Item(9, 6).Value2 = 1
$ws.Cells.Item(9, 7).Value2 = 20.82581466666667
$ws.Cells.Item(9, 8).Value2 = 62.47744400000001
$ws.Cells.Item(9, 9).Value2 = 0.9748457050079848
$ws.Cells.Item(9, 10).Value2 = 0.9748457050079848
$ws.Cells.Item(9, 11).Value2 = 3
$ws.Cells.Item(9, 12).Value2 = 1
$ws.Cells.Item(9, 13).Value2 = 0.2815866666666667
$ws.Cells.Item(9, 14).Value2 = 0.8447600000000001
$ws.Cells.Item(9, 15).Value2 = 0.05434147372589489
$ws.Cells.Item(9, 16).Value2 = 0.05434147372589488
$ws.Cells.Item(9, 17).Value2 = 5.864271732604446
$ws.Cells.Item(9, 18).Value2 = 52.77844559344001
$ws.Cells.Item(9, 19).Value2 = 0.05297455226549289
$ws.Cells.Item(9, 20).Value2 = 0.05297455226549288

# Row 10
$ws.Cells.Item(10, 1).Value2 = "MuSCs"
$ws.Cells.Item(10, 2).Value2 = "Clec11a"
$ws.Cells.Item(10, 3).Value2 = "Itga10"
$ws.Cells.Item(10, 4).Value2 = "ECs"
$ws.Cells.Item(10, 5).Value2 = 3
$ws.Cells.Item(10, 6).Value2 = 1
$ws.Cells.Item(10, 7).Value2 = 0.200477
$ws.Cells.Item(10, 8).Value2 = 0.601431
$ws.Cells.Item(10, 9).Value2 = 0.009384225564807953
$ws.Cells.Item(10, 10).Value2 = 0.009384225564807953
$ws.Cells.Item(10, 11).Value2 = 3
$ws.Cells.Item(10, 12).Value2 = 1
$ws.Cells.Item(10, 13).Value2 = 0.970007
$ws.Cells.Item(10, 14).Value2 = 2.910021
$ws.Cells.Item(10, 15).Value2 = 0.1871949781160357
$ws.Cells.Item(10, 16).Value2 = 0.1871949781160357
$ws.Cells.Item(10, 17).Value2 = 0.194464093339
$ws.Cells.Item(10, 18).Value2 = 1.750176840051
$ws.Cells.Item(10, 19).Value2 = 0.001756679899240168
$ws.Cells.Item(10, 20).Value2 = 0.001756679899240168

# Row 11
$ws.Cells.Item(11, 1).Value2 = "MuSCs"
$ws.Cells.Item(11, 2).Value2 = "Clec11a"
$ws.Cells.Item(11, 3).Value2 = "Itga10"
$ws.Cells.Item(11, 4).Value2 = "FAPs"
$ws.Cells.Item(11, 5).Value2 = 3
$ws.Cells.Item(11, 6).Value2 = 1
$ws.Cells.Item(11, 7).Value2 = 0.200477
$ws.Cells.Item(11, 8).Value2 = 0.601431
$ws.Cells.Item(11, 9).Value2 = 0.009384225564807953
$ws.Cells.Item(11, 10).Value2 = 0.009384225564807953
$ws.Cells.Item(11, 11).Value2 = 3
$ws.Cells.Item(11, 12).Value2 = 1
$ws.Cells.Item(11, 13).Value2 = 2.313512
$ws.Cells.Item(11, 14).Value2 = 6.940536
$ws.Cells.Item(11, 15).Value2 = 0.4464687659070358
$ws.Cells.Item(11, 16).Value2 = 0.4464687659070357
$ws.Cells.Item(11, 17).Value2 = 0.463805945224
$ws.Cells.Item(11, 18).Value2 = 4.174253507016
$ws.Cells.Item(11, 19).Value2 = 0.004189763606913063
$ws.Cells.Item(11, 20).Value2 = 0.004189763606913062

# Row 12
$ws.Cells.Item(12, 1).Value2 = "MuSCs"
$ws.Cells.Item(12, 2).Value2 = "Clec11a"
$ws.Cells.Item(12, 3).Value2 = "Itga10"
$ws.Cells.Item(12, 4).Value2 = "MuSCs"
$ws.Cells.Item(12, 5).Value2 = 3
$ws.Cells.Item(12, 6).Value2 = 1
$ws.Cells.Item(12, 7).Value2 = 0.200477
$ws.Cells.Item(12, 8).Value2 = 0.601431
$ws.Cells.Item(12, 9).Value2 = 0.009384225564807953
$ws.Cells.Item(12, 10).Value2 = 0.009384225564807953
$ws.Cells.Item(12, 11).Value2 = 3
$ws.Cells.Item(12, 12).Value2 = 1
$ws.Cells.Item(12, 13).Value2 = 1.616694666666667
$ws.Cells.Item(12, 14).Value2 = 4.850084
$ws.Cells.Item(12, 15).Value2 = 0.3119947822510337
$ws.Cells.Item(12, 16).Value2 = 0.3119947822510336
$ws.Cells.Item(12, 17).Value2 = 0.3241100966893333
$ws.Cells.Item(12, 18).Value2 = 2.916990870204
$ws.Cells.Item(12, 19).Value2 = 0.002927829411686841
$ws.Cells.Item(12, 20).Value2 = 0.00292782941168684

# Row 13
$ws.Cells.Item(13, 1).Value2 = "MuSCs"
$ws.Cells.Item(13, 2).Value2 = "Clec11a"
$ws.Cells.Item(13, 3).Value2 = "Itga10"
$ws.Cells.Item(13, 4).Value2 = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value2 = 3
$ws.Cells.Item(13, 6).Value2 = 1
$ws.Cells.Item(13, 7).Value2 = 0.200477
$ws.Cells.Item(13, 8).Value2 = 0.601431
$ws.Cells.Item(13, 9).Value2 = 0.009384225564807953
$ws.Cells.Item(13, 10).Value2 = 0.009384225564807953
$ws.Cells.Item(13, 11).Value2 = 3
$ws.Cells.Item(13, 12).Value2 = 1
$ws.Cells.Item(13, 13).Value2 = 0.2815866666666667
$ws.Cells.Item(13, 14).Value2 = 0.8447600000000001
$ws.Cells.Item(13, 15).Value2 = 0.05434147372589489
$ws.Cells.Item(13, 16).Value2 = 0.05434147372589488
$ws.Cells.Item(13, 17).Value2 = 0.05645165017333335
$ws.Cells.Item(13, 18).Value2 = 0.50806485156
$ws.Cells.Item(13, 19).Value2 = 0.0005099526469678825
$ws.Cells.Item(13, 20).Value2 = 0.0005099526469678824

# Row 14
$ws.Cells.Item(14, 1).Value2 = "Resolving-Mac"
$ws.Cells.Item(14, 2).Value2 = "Clec11a"
$ws.Cells.Item(14, 3).Value2 = "Itga10"
$ws.Cells.Item(14, 4).Value2 = "ECs"
$ws.Cells.Item(14, 5).Value2 = 2
$ws.Cells.Item(14, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(14, 7).Value2 = 0.183981
$ws.Cells.Item(14, 8).Value2 = 0.5519430000000001
$ws.Cells.Item(14, 9).Value2 = 0.008612056264005009
$ws.Cells.Item(14, 10).Value2 = 0.008612056264005009
$ws.Cells.Item(14, 11).Value2 = 3
$ws.Cells.Item(14, 12).Value2 = 1
$ws.Cells.Item(14, 13).Value2 = 0.970007
$ws.Cells.Item(14, 14).Value2 = 2.910021
$ws.Cells.Item(14, 15).Value2 = 0.1871949781160357
$ws.Cells.Item(14, 16).Value2 = 0.1871949781160357
$ws.Cells.Item(14, 17).Value2 = 0.178462857867
$ws.Cells.Item(14, 18).Value2 = 1.606165720803
$ws.Cells.Item(14, 19).Value2 = 0.001612133683874486
$ws.Cells.Item(14, 20).Value2 = 0.001612133683874486

# Row 15
$ws.Cells.Item(15, 1).Value2 = "Resolving-Mac"
$ws.Cells.Item(15, 2).Value2 = "Clec11a"
$ws.Cells.Item(15, 3).Value2 = "Itga10"
$ws.Cells.Item(15, 4).Value2 = "FAPs"
$ws.Cells.Item(15, 5).Value2 = 2
$ws.Cells.Item(15, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(15, 7).Value2 = 0.183981
$ws.Cells.Item(15, 8).Value2 = 0.5519430000000001
$ws.Cells.Item(15, 9).Value2 = 0.008612056264005009
$ws.Cells.Item(15, 10).Value2 = 0.008612056264005009
$ws.Cells.Item(15, 11).Value2 = 3
$ws.Cells.Item(15, 12).Value2 = 1
$ws.Cells.Item(15, 13).Value2 = 2.313512
$ws.Cells.Item(15, 14).Value2 = 6.940536
$ws.Cells.Item(15, 15).Value2 = 0.4464687659070358
$ws.Cells.Item(15, 16).Value2 = 0.4464687659070357
$ws.Cells.Item(15, 17).Value2 = 0.4256422512720001
$ws.Cells.Item(15, 18).Value2 = 3.830780261448
$ws.Cells.Item(15, 19).Value2 = 0.003845014132112274
$ws.Cells.Item(15, 20).Value2 = 0.003845014132112273

# Row 16
$ws.Cells.Item(16, 1).Value2 = "Resolving-Mac"
$ws.Cells.Item(16, 2).Value2 = "Clec11a"
$ws.Cells.Item(16, 3).Value2 = "Itga10"
$ws.Cells.Item(16, 4).Value2 = "MuSCs"
$ws.Cells.Item(16, 5).Value2 = 2
$ws.Cells.Item(16, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(16, 7).Value2 = 0.183981
$ws.Cells.Item(16, 8).Value2 = 0.5519430000000001
$ws.Cells.Item(16, 9).Value2 = 0.008612056264005009
$ws.Cells.Item(16, 10).Value2 = 0.008612056264005009
$ws.Cells.Item(16, 11).Value2 = 3
$ws.Cells.Item(16, 12).Value2 = 1
$ws.Cells.Item(16, 13).Value2 = 1.616694666666667
$ws.Cells.Item(16, 14).Value2 = 4.850084
$ws.Cells.Item(16, 15).Value2 = 0.3119947822510337
$ws.Cells.Item(16, 16).Value2 = 0.3119947822510336
$ws.Cells.Item(16, 17).Value2 = 0.2974411014680001
$ws.Cells.Item(16, 18).Value2 = 2.676969913212
$ws.Cells.Item(16, 19).Value2 = 0.002686916618821893
$ws.Cells.Item(16, 20).Value2 = 0.002686916618821892

# Row 17
$ws.Cells.Item(17, 1).Value2 = "Resolving-Mac"
$ws.Cells.Item(17, 2).Value2 = "Clec11a"
$ws.Cells.Item(17, 3).Value2 = "Itga10"
$ws.Cells.Item(17, 4).Value2 = "Resolving-Mac"
$ws.Cells.Item(17, 5).Value2 = 2
$ws.Cells.Item(17, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(17, 7).Value2 = 0.183981
$ws.Cells.Item(17, 8).Value2 = 0.5519430000000001
$ws.Cells.Item(17, 9).Value2 = 0.008612056264005009
$ws.Cells.Item(17, 10).Value2 = 0.008612056264005009
$ws.Cells.Item(17, 11).Value2 = 3
$ws.Cells.Item(17, 12).Value2 = 1
$ws.Cells.Item(17, 13).Value2 = 0.2815866666666667
$ws.Cells.Item(17, 14).Value2 = 0.8447600000000001
$ws.Cells.Item(17, 15).Value2 = 0.05434147372589489
$ws.Cells.Item(17, 16).Value2 = 0.05434147372589488
$ws.Cells.Item(17, 17).Value2 = 0.05180659652000002
$ws.Cells.Item(17, 18).Value2 = 0.4662593686800001
$ws.Cells.Item(17, 19).Value2 = 0.0004679918291963567
$ws.Cells.Item(17, 20).Value2 = 0.0004679918291963566

